$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5,D9,D10,D11,D15,D16,D19,D21,D23,D24,D25,D28,D30,D32,D39,D42,D45,D46').NumberFormat = '@'

$ws.Range('D2').Value = '26.852.11'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '1.641.29'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('E4').Value = '  -0.60%  '
$ws.Range('D5').Value = '216.23'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('E6').Value = '  +1.74%  '
$ws.Range('E7').Value = '  -0.58%  '
$ws.Range('E8').Value = '  +1.77%  '
$ws.Range('D9').Value = '0.0620'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').Value = '19.84'
$ws.Range('E10').Value = '  +4.61%  '
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').Value = '1.871.21'
$ws.Range('D13').Value = '1.643.16'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').Value = '0.527'
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '66.39'
$ws.Range('E16').Value = '  +3.85%  '
$ws.Range('D17').Value = '26.853.09'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +1.46%  '
$ws.Range('D19').Value = '218.34'
$ws.Range('E19').Value = '  +3.44%  '
$ws.Range('E20').Value = '  -0.61%  '
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  +8.25%  '
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('D23').Value = '2.42'
$ws.Range('E23').Value = '  +4.15%  '
$ws.Range('D24').Value = '9.16'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = '145.97'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('E27').Value = '  +5.15%  '
$ws.Range('D28').Value = '0.119'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('E29').Value = '  +1.84%  '
$ws.Range('D30').Value = '0.0510'
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('E34').Value = '  +2.90%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '1.236.91'
$ws.Range('E37').Value = '  +1.15%  '
$ws.Range('E38').Value = '  +3.74%  '
$ws.Range('D39').Value = '0.833'
$ws.Range('E39').Value = '  +4.59%  '
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D42').Value = '5.36'
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('D43').Value = '1.782.34'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('E44').Value = '  -2.95%  '
$ws.Range('D45').Value = '60.86'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('D46').Value = '91.50'
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +18.51%  '
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('E50').Value = '  +1.74%  '
$ws.Range('E51').Value = '  +1.91%  '
